# dev-1.0.0 : Update Templete dan Index Sales
#
# Fill in the monthly (E:P) figures for every detail row (2-22) of the
# "Sales" sheet with 1, carrying over the already-highlighted input style
# (fillId=5 / numFmt "#,##0") that a handful of cells (E2, E4, E6, ...)
# already used, and update the sheet's active selection to reflect the
# range that was just filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")
$ws.Activate()

# E2 already carries the exact "highlighted input" number format (blue
# fill, thousands separator, thin grid border) that every other monthly
# cell in the table should end up with, so use it as the format template.
$formatSource = $ws.Range("E2")

$firstRow = 2
$lastRow = 22

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowRange = $ws.Range("E" + $r + ":P" + $r)
    $rowRange.Value = 1
}

# Copy E2's formatting (fill/border/number-format) across the whole
# E2:P22 block so every cell matches the existing "input" look.
$formatSource.Copy() | Out-Null
$ws.Range("E2:P22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Recalculate the Q/R/S rollup formulas now that the inputs changed.
$excel.Calculate() | Out-Null

# Reflect the edited block in the sheet's view/selection.
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E2:P22").Select() | Out-Null
